$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the expense amount columns (D15:E22) with a yellow fill -
# this mints the new solid-yellow fill + cellXf used throughout the block.
$ws.Range("D15:E22").Interior.Color = 65535

# --- Row 21: new "Other" expense line (reuses the existing "Gasoline for
# system" shared string, same as row 20) ---
$ws.Range("B21").Value = "Gasoline for system"
$ws.Range("I21").Value = 22

# --- Row 22: new "Super Shuttle from SFO" expense line ---
# Copy the date style (s="7", numFmtId 14) from a cell above before writing
# the serial date value so A22 matches A15:A20 formatting.
$ws.Range("A17").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 42390
$ws.Range("B22").Value = "Super Shuttle from SFO"
$ws.Range("E22").Value = 41.3

# --- Totals ---
$ws.Range("D24").Formula = "=SUM(D15:D16)"
$ws.Range("E26").Formula = "=SUM(D15:E22)"

# --- Selection state matches the saved view in the target workbook ---
[void]$ws.Range("D15:E22").Select()
